$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for 2000年, 2002年, 2005年, 2007年 (old rows 2-5).
# This shifts 2010年/2012年/2015年/2017年 (old rows 6-9) up to rows 2-5.
$ws.Rows("2:5").Delete()

# Append a new row for 2020年 as row 6, reusing the year-label formatting
# (bold font, thin box border, centered/top aligned) from the row above.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "2020年"
$ws.Range("B6").Value = -157992.322467643
$ws.Range("C6").Value = 5234595.238605
$ws.Range("D6").Value = 72478841.18828399
$ws.Range("F6").Value = 2707176325.74967
$ws.Range("I6").Value = 820471959.731071
$ws.Range("J6").Value = -25605149.2034971
$ws.Range("K6").Value = 1785150.09795483
$ws.Range("L6").Value = -3571111.18213603
$ws.Range("M6").Value = 4335213941.04906
$ws.Range("O6").Value = 462218.578754962
$ws.Range("P6").Value = 38119950.9600789
$ws.Range("R6").Value = 2140696.25007459
$ws.Range("S6").Value = 8525640.232746361
